$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value2 = 425.1
$ws.Range("I2").Value2 = 444.55554
$ws.Range("J2").Value2 = 250
$ws.Range("K2").Value2 = 444.55554
$ws.Range("L2").Value2 = 250
$ws.Range("M2").Value2 = -331.55554
$ws.Range("N2").Value2 = -476

$ws.Range("H8").Value2 = 43
$ws.Range("I8").Value2 = 43
$ws.Range("K8").Value2 = 129
$ws.Range("M8").Value2 = 10

$ws.Range("H41").Value2 = 285.7143
$ws.Range("I41").Value2 = 100
$ws.Range("J41").Value2 = 750
$ws.Range("K41").Value2 = 100
$ws.Range("L41").Value2 = 750
$ws.Range("M41").Value2 = 340
$ws.Range("N41").Value2 = -1630

$ws.Range("H43").Value2 = 2727.8
$ws.Range("I43").Value2 = 2400
$ws.Range("J43").Value2 = 2809.75
$ws.Range("K43").Value2 = 2400
$ws.Range("L43").Value2 = 2809.75
$ws.Range("M43").Value2 = -2331
$ws.Range("N43").Value2 = -2947.75

$ws.Range("H53").Value2 = 141.51613
$ws.Range("I53").Value2 = 104.8
$ws.Range("J53").Value2 = 175.9375
$ws.Range("K53").Value2 = 104.8
$ws.Range("L53").Value2 = 175.9375
$ws.Range("M53").Value2 = 532.2
$ws.Range("N53").Value2 = -1449.9375

$ws.Range("H62").Value2 = 3081
$ws.Range("I62").Value2 = 3101.25
$ws.Range("J62").Value2 = 3000
$ws.Range("K62").Value2 = 3101.25
$ws.Range("L62").Value2 = 3000
$ws.Range("M62").Value2 = -2477.25
$ws.Range("N62").Value2 = -4248

$ws.Range("H65").Value2 = 3081
$ws.Range("I65").Value2 = 3101.25
$ws.Range("J65").Value2 = 3000
$ws.Range("K65").Value2 = 15506.25
$ws.Range("L65").Value2 = 15000
$ws.Range("M65").Value2 = -12386.25
$ws.Range("N65").Value2 = -21240

$ws.Range("H129").Value2 = 1097.44
$ws.Range("I129").Value2 = 333.33334
$ws.Range("J129").Value2 = 1146.2128
$ws.Range("K129").Value2 = 1000.00002
$ws.Range("L129").Value2 = 3438.6384
$ws.Range("M129").Value2 = 3999.99998
$ws.Range("N129").Value2 = -13438.6384

$ws.Range("H132").Value2 = 1940.7576
$ws.Range("I132").Value2 = 2090.8928
$ws.Range("J132").Value2 = 1100
$ws.Range("K132").Value2 = 6272.678400000001
$ws.Range("L132").Value2 = 3300
$ws.Range("M132").Value2 = -3742.678400000001
$ws.Range("N132").Value2 = -8360

$ws.Range("H136").Value2 = 49000
$ws.Range("J136").Value2 = 49000
$ws.Range("L136").Value2 = 49000
$ws.Range("N136").Value2 = -59200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 23025.121
$ws.Range("I32").Value2 = 26203.111
$ws.Range("J32").Value2 = 12024.385
$ws.Range("K32").Value2 = 26203.111
$ws.Range("L32").Value2 = 12024.385
$ws.Range("M32").Value2 = -25916.111
$ws.Range("N32").Value2 = -12598.385

$ws.Range("H45").Value2 = 1043.8704
$ws.Range("I45").Value2 = 921.3022999999999
$ws.Range("K45").Value2 = 921.3022999999999
$ws.Range("M45").Value2 = -544.3022999999999

$ws.Range("H118").Value2 = 0
$ws.Range("J118").Value2 = 0
$ws.Range("L118").Value2 = 0
$ws.Range("N118").ClearContents()

$ws.Range("H122").Value2 = 1534.5834
$ws.Range("I122").Value2 = 1021.4667
$ws.Range("J122").Value2 = 2389.7778
$ws.Range("K122").Value2 = 3064.4001
$ws.Range("L122").Value2 = 7169.3334
$ws.Range("M122").Value2 = -614.4000999999998
$ws.Range("N122").Value2 = -12069.3334

$ws.Range("H132").Value2 = 2613.96
$ws.Range("I132").Value2 = 2028.2222
$ws.Range("J132").Value2 = 4120.143
$ws.Range("K132").Value2 = 6084.6666
$ws.Range("L132").Value2 = 12360.429
$ws.Range("M132").Value2 = -3554.6666
$ws.Range("N132").Value2 = -17420.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value2 = 1400.6666
$ws.Range("I107").Value2 = 1400.6666
$ws.Range("J107").Value2 = 0
$ws.Range("K107").Value2 = 1400.6666
$ws.Range("L107").Value2 = 0
$ws.Range("M107").Value2 = 519.3334
$ws.Range("N107").ClearContents()

$ws.Range("H134").Value2 = 2702.7576
$ws.Range("I134").Value2 = 2440.4075
$ws.Range("K134").Value2 = 7321.2225
$ws.Range("M134").Value2 = -4786.2225

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value2 = 14380.429
$ws.Range("J122").Value2 = 25112.8
$ws.Range("L122").Value2 = 75338.39999999999
$ws.Range("N122").Value2 = -80238.39999999999

$ws.Range("H139").Value2 = 0
$ws.Range("J139").Value2 = 0
$ws.Range("L139").Value2 = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value2 = 21306.105
$ws.Range("J131").Value2 = 24850.426
$ws.Range("L131").Value2 = 74551.27799999999
$ws.Range("N131").Value2 = -84631.27799999999

$ws.Range("H137").Value2 = 18671.912
$ws.Range("I137").Value2 = 1500.6
$ws.Range("J137").Value2 = 32228.21
$ws.Range("K137").Value2 = 4501.799999999999
$ws.Range("L137").Value2 = 96684.63
$ws.Range("M137").Value2 = 598.2000000000007
$ws.Range("N137").Value2 = -106884.63

$ws.Range("H140").Value2 = 2378.4773
$ws.Range("I140").Value2 = 1848.7931
$ws.Range("J140").Value2 = 3402.5334
$ws.Range("K140").Value2 = 5546.379300000001
$ws.Range("L140").Value2 = 10207.6002
$ws.Range("M140").Value2 = -366.3793000000005
$ws.Range("N140").Value2 = -20567.6002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value2 = 2780.6667
$ws.Range("I113").Value2 = 3351.8333
$ws.Range("K113").Value2 = 3351.8333
$ws.Range("M113").Value2 = -1181.8333

$ws.Range("H132").Value2 = 3564.9375
$ws.Range("I132").Value2 = 3234.6924
$ws.Range("J132").Value2 = 4996
$ws.Range("K132").Value2 = 9704.0772
$ws.Range("L132").Value2 = 14988
$ws.Range("M132").Value2 = -7174.0772
$ws.Range("N132").Value2 = -20048

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 5706.4707
$ws.Range("I7").Value2 = 5257.5
$ws.Range("J7").Value2 = 7801.6665
$ws.Range("K7").Value2 = 5257.5
$ws.Range("L7").Value2 = 7801.6665
$ws.Range("M7").Value2 = -5145.5
$ws.Range("N7").Value2 = -8025.6665

$ws.Range("H61").Value2 = 602580.9399999999
$ws.Range("I61").Value2 = 16158.066
$ws.Range("J61").Value2 = 5000752.5
$ws.Range("K61").Value2 = 16158.066
$ws.Range("L61").Value2 = 5000752.5
$ws.Range("M61").Value2 = -15956.066
$ws.Range("N61").Value2 = -5001156.5

$ws.Range("H113").Value2 = 602580.9399999999
$ws.Range("I113").Value2 = 16158.066
$ws.Range("J113").Value2 = 5000752.5
$ws.Range("K113").Value2 = 16158.066
$ws.Range("L113").Value2 = 5000752.5
$ws.Range("M113").Value2 = -13988.066
$ws.Range("N113").Value2 = -5005092.5

$ws.Range("H122").Value2 = 6255.6763
$ws.Range("I122").Value2 = 5960.289
$ws.Range("J122").Value2 = 6833.609
$ws.Range("K122").Value2 = 17880.867
$ws.Range("L122").Value2 = 20500.827
$ws.Range("M122").Value2 = -15430.867
$ws.Range("N122").Value2 = -25400.827

$ws.Range("H126").Value2 = 5706.4707
$ws.Range("I126").Value2 = 5257.5
$ws.Range("J126").Value2 = 7801.6665
$ws.Range("K126").Value2 = 15772.5
$ws.Range("L126").Value2 = 23404.9995
$ws.Range("M126").Value2 = -13302.5
$ws.Range("N126").Value2 = -28344.9995

$ws.Range("H136").Value2 = 3857.8276
$ws.Range("I136").Value2 = 2102.5715
$ws.Range("J136").Value2 = 6528.8696
$ws.Range("K136").Value2 = 6307.7145
$ws.Range("L136").Value2 = 19586.6088
$ws.Range("M136").Value2 = -3757.7145
$ws.Range("N136").Value2 = -24686.6088

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value2 = 8746.923000000001
$ws.Range("I113").Value2 = 33667
$ws.Range("J113").Value2 = 1270.9
$ws.Range("K113").Value2 = 101001
$ws.Range("L113").Value2 = 3812.7
$ws.Range("M113").Value2 = -98831
$ws.Range("N113").Value2 = -8152.700000000001

$ws.Range("H115").Value2 = 50000
$ws.Range("J115").Value2 = 50000
$ws.Range("L115").Value2 = 50000
$ws.Range("N115").Value2 = -53134

$ws.Range("H122").Value2 = 3289.1904
$ws.Range("I122").Value2 = 2020.5
$ws.Range("K122").Value2 = 6061.5
$ws.Range("M122").Value2 = -3611.5
